$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells, copy style from existing header cell (H1) for formatting consistency
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats

# Data values for I and J columns (rows 2-17)
$values = @(
    @(7, 7),
    @(10, 10),
    @(8, 8),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(6, 7),
    @(6, 7),
    @(6, 7),
    @(9, 9),
    @(9, 9),
    @(6, 7),
    @(8, 8),
    @(9, 9),
    @(4, 4),
    @(4, 4)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
